# "Two Trigger integrated with front-end and working superbly fine"
# 1) Email Validation Trigger Done
# 2) Password Change Trigger Done
# Front-end improved and made more dynamic.
#
# This script reproduces, via Excel COM-interop, the edits made to the
# "DD" (database design) worksheet: two existing trigger name cells get
# left-aligned, and a new worked example / truth-table for the
# "Table : Order_detail" cart/checkout/order flow is added below the
# existing Order_detail table (rows 57-75, columns F/G), along with the
# six new shared strings that content requires.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DD")

# --- Left-align the two already-existing trigger name cells (the new
#     "Email validation" / "Password change" triggers being called out) ---
$ws.Range("F19").HorizontalAlignment = -4131   # xlLeft - user_singup_validation
$ws.Range("F21").HorizontalAlignment = -4131   # xlLeft - change_password
$ws.Range("F23").HorizontalAlignment = -4131   # xlLeft - price_cannot_be_negative

# --- New "cart" mini table (headers + a 1-row qty truth table) ---
$ws.Range("F57").Value = "cart "
$ws.Range("F58").Value = "productid"
$ws.Range("G58").Value = "qty"

$ws.Range("F59").Value = 1
$ws.Range("G59").Value = 2
$ws.Range("G59").HorizontalAlignment = -4131   # xlLeft

$ws.Range("F60").Value = 2
$ws.Range("G60").Value = 2
$ws.Range("G60").HorizontalAlignment = -4131   # xlLeft

$ws.Range("F61").Value = 3
$ws.Range("G61").Value = 6
$ws.Range("G61").HorizontalAlignment = -4131   # xlLeft

# --- "checkout" label next to the Order_detail table ---
$ws.Range("F63").Value = "checkout"

# --- "order" section with a small id/qty example ---
$ws.Range("F67").Value = "order "

$ws.Range("E69").Value = "id"
$ws.Range("F69").Value = 1
$ws.Range("G69").Value = 500

# --- "orderdetail" section with its own worked example ---
$ws.Range("F71").Value = "orderdetail"
$ws.Range("F72").Value = 1
$ws.Range("F73").Value = 1
$ws.Range("F74").Value = 2
$ws.Range("F75").Value = 500

# --- Update the on-screen selection/scroll to match the author's saved
#     view (scrolled right to the new content, cursor on F11) ---
$ws.Range("F11").Select()
